$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# About sheet
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Insert a new row at row 11. Excel shifts everything below down by one and
# automatically repoints any formula that referenced About!$A$16 (the old
# conversion-factor cell) to About!$A$17 -- including the cross-sheet
# formulas on SourceData.
$about.Rows.Item(11).Insert()

# Row 10 text is reworded (split across two sentences now).
$about.Range("A10").Value = "When considering the Social Cost of Carbon, meant to capture the long-term economic damage caused by one"

# New row 11 holds the second half of the reworded note.
$about.Range("A11").Value = "ton of carbon dioxide emitted, the U.S. government typically uses the figures based on"

# Row 16 (formerly row 15) conversion-factor note: 2017 -> 2012 dollars.
$about.Range("A16").Value = "We adjust 2007 dollars to 2012 dollars using the following conversion factor:"

# Row 17 (formerly row 16) used to hold the formula =1.109*1.068; it is now a
# plain literal value (no 1.068 CPI adjustment applied any more).
$about.Range("A17").Value = 1.109

# ---------------------------------------------------------------------------
# SCoC sheet
# ---------------------------------------------------------------------------
$scoc = $wb.Worksheets.Item("SCoC")
$scoc.Range("B1").Value = "Social Cost of Carbon ($/g CO2e)"

$wb.Application.Calculate()
